# Recalibrated SABR parameters (alpha, rho, nu) — updates the 3x5 grid
# (rows 1Y/5Y/10Y x tenors 1Y/2Y/3Y/5Y/10Y) on each of the three sheets.

$wb = $excel.ActiveWorkbook

function Set-Grid($ws, $values) {
    for ($r = 0; $r -lt $values.Length; $r++) {
        $row = $values[$r]
        for ($c = 0; $c -lt $row.Length; $c++) {
            $ws.Cells.Item($r + 2, $c + 2).Value = $row[$c]
        }
    }
}

# --- sabr_alpha ---
$ws = $wb.Worksheets.Item("sabr_alpha")
Set-Grid $ws @(
    @(0.1391035778845959, 0.1846425592601795, 0.1967712904214697, 0.1779273404887576, 0.172890201220095),
    @(0.1664545218569979, 0.1983862754120166, 0.1999508905666859, 0.1880623432715516, 0.171354499666421),
    @(0.1783943026753877, 0.1932472883083041, 0.1968208170551125, 0.1929909483890378, 0.1789472682163097)
)

# --- sabr_rho ---
$ws = $wb.Worksheets.Item("sabr_rho")
Set-Grid $ws @(
    @(-0.6330062548889566, -0.5250855973782115, -0.4826171172249953, -0.4136425277006885, -0.3097880589022577),
    @(-0.5851338342383721, -0.5438848697963986, -0.5179597068906967, -0.4957420723203845, -0.3971619656711662),
    @(-0.5488092695261677, -0.540522024976866,  -0.5276708274378461, -0.540684142814833,  -0.5127356878951527)
)

# --- sabr_nu ---
$ws = $wb.Worksheets.Item("sabr_nu")
Set-Grid $ws @(
    @(2.048783209110644, 1.677390418935748,  1.438553967961835,  1.066043910350141,  0.7604156618523644),
    @(1.339865494949765, 1.062308192887452,  0.9467116104193644, 0.6835419435497165, 0.5433373638486673),
    @(1.010697794727172, 0.9245785503057192, 0.8608700206026637, 0.7236644034101184, 0.6015151318353243)
)
